# Separate utils from main package
# Re-run of the logging process: refresh the start/end execution timestamps
# (columns E/F) and the computed duration text (column G) for log rows 2-10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 2;  E = 45130.66270084777;  F = 45130.66270084898;  G = $null },
    @{ Row = 3;  E = 45130.66270534682;  F = 45130.66270877632;  G = "0d, 0hr, 0min, 0.296sec " },
    @{ Row = 4;  E = 45130.66270921857;  F = 45130.66271368905;  G = "0d, 0hr, 0min, 0.386sec " },
    @{ Row = 5;  E = 45130.66271415927; F = 45130.66271602872;  G = "0d, 0hr, 0min, 0.161sec " },
    @{ Row = 6;  E = 45130.66271637068;  F = 45130.66271653711;  G = "0d, 0hr, 0min, 0.014sec " },
    @{ Row = 7;  E = 45130.66271706452;  F = 45130.66272009206;  G = "0d, 0hr, 0min, 0.261sec " },
    @{ Row = 8;  E = 45130.66272048949;  F = 45130.66272637997;  G = "0d, 0hr, 0min, 0.508sec " },
    @{ Row = 9;  E = 45130.66272676507;  F = 45130.66273586931;  G = "0d, 0hr, 0min, 0.786sec " },
    @{ Row = 10; E = 45130.66273630778; F = 45130.66274638397;  G = "0d, 0hr, 0min, 0.87sec " }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 5).Value = $r.E
    $ws.Cells.Item($r.Row, 6).Value = $r.F
    if ($r.G -ne $null) {
        $ws.Cells.Item($r.Row, 7).Value = $r.G
    }
}
